$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = 6.598000000000001
    12 = 5.492999999999999
    15 = 5.295
    27 = 6.709000000000001
    28 = 6.726000000000001
    31 = 6.796000000000001
    32 = 7.341999999999999
    36 = 8.704000000000001
    38 = 6.052
    46 = 6.477000000000001
    54 = 5.05
    55 = 4.459000000000001
    56 = 4.700000000000001
    67 = 4.988999999999999
    69 = 5.002000000000001
    72 = 6.347
    73 = 6.600999999999999
    83 = 5.512
    86 = 5.146000000000001
    91 = 6.109
    93 = 5.380999999999999
    99 = 5.217000000000001
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
